$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '26.263.77'
$ws.Range("E2").Value = '  +0.66%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.665.00'
$ws.Range("E3").Value = '  +0.62%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.010'
$ws.Range("E4").Value = '  +0.87%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '218.69'
$ws.Range("E5").Value = '  +0.22%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.5302'
$ws.Range("E6").Value = '  +0.46%  '
$ws.Range("E7").Value = '  +0.86%  '
$ws.Range("E8").Value = '  +1.11%  '
$ws.Range("E9").Value = '  +0.44%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '20.56'
$ws.Range("E10").Value = '  +0.54%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07849'
$ws.Range("E11").Value = '  +0.82%  '
$ws.Range("E12").Value = '  +1.47%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.663.91'
$ws.Range("E13").Value = '  +0.54%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '1.891.72'
$ws.Range("E14").Value = '  +0.51%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.5527'
$ws.Range("E15").Value = '  +0.74%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0₅8175'
$ws.Range("E16").Value = '  +0.10%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '65.61'
$ws.Range("E17").Value = '  +0.31%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '1.010'
$ws.Range("E18").Value = '  +0.81%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '4.670'
$ws.Range("E19").Value = '  +2.47%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '193.51'
$ws.Range("E20").Value = '  +0.49%  '
$ws.Range("E21").Value = '  +1.53%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.044'
$ws.Range("E22").Value = '  +0.22%  '
$ws.Range("E23").Value = '  +0.80%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '144.60'
$ws.Range("E24").Value = '  +2.35%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.1226'
$ws.Range("E25").Value = '  -1.78%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '7.215'
$ws.Range("E26").Value = '  -0.74%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '16.09'
$ws.Range("E27").Value = '  -0.79%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.483'
$ws.Range("E28").Value = '  +3.43%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.06001'
$ws.Range("E29").Value = '  +1.39%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.284'
$ws.Range("E30").Value = '  +0.30%  '
$ws.Range("E31").Value = '  +1.71%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.287'
$ws.Range("E32").Value = '  +1.27%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.614'
$ws.Range("E33").Value = '  +2.90%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.9611'
$ws.Range("E34").Value = '  +1.31%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.826'
$ws.Range("E35").Value = '  +0.76%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.425'
$ws.Range("E36").Value = '  +0.66%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.5811'
$ws.Range("E37").Value = '  +2.93%  '
$ws.Range("E38").Value = '  -0.29%  '
$ws.Range("B39").Value = 'TrustWalletToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.8657'
$ws.Range("E39").Value = '  +2.21%  '
$ws.Range("B40").Value = 'FraxShare'
$ws.Range("C40").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '5.871'
$ws.Range("E40").Value = '  +0.96%  '
$ws.Range("B41").Value = 'Maker'
$ws.Range("C41").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.050.05'
$ws.Range("E41").Value = '  +2.42%  '
$ws.Range("B42").Value = 'PaxDollar'
$ws.Range("C42").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.010'
$ws.Range("E42").Value = '  +0.83%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '104.30'
$ws.Range("E43").Value = '  +2.21%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.804.35'
$ws.Range("E44").Value = '  +0.38%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '57.43'
$ws.Range("E45").Value = '  +0.56%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.016'
$ws.Range("E46").Value = '  +0.76%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0₈105'
$ws.Range("E47").Value = '  -5.36%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.4383'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '7.994'
$ws.Range("E49").Value = '  +2.54%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.05166'
$ws.Range("E50").Value = '  +0.30%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.439'
$ws.Range("E51").Value = '  -2.30%  '
